# Auto commit at 2025-11-15  8:35:14.12
# Append two new daily rows (150 & 151) of hourly charging data for
# 2025-11-14 ("45975") to Sheet1, one row per station, mirroring the
# layout/styles of the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Station labels reuse the existing shared-string text used by every other
# row in column B.
$station1 = "四方坪站充电量(kw)"
$station2 = "高岭站充电量(kw)"

# New date serial value (2025-11-14) shared by both rows, same as the
# canonical Excel date serial stored in column A for prior rows.
$dateValue = 45975

$row150 = @(603.87599999999998,982.20799999999997,598.04199999999992,288.44600000000003,377.75300000000004,704.98699999999997,540.93799999999999,185.40600000000001,197.84800000000001,156.93299999999999,121.42,222.60699999999997,767.0200000000001,1226.1940000000002,587.59999999999991,471.44800000000004,345.44400000000007,225.9,280.17599999999999,68.02000000000001,66.44,15.9,95.789999999999992,49.595999999999997)

$row151 = @(489.9729999999999,329.01599999999996,23.039000000000001,116.78199999999998,55.382999999999996,292.858,139.69599999999997,62.625,393.28800000000001,213.70599999999999,164.36,344.11599999999999,419.209,512.44399999999996,146.80300000000003,299.72300000000001,131.14400000000001,165.07300000000001,38.938000000000002,111.524,44.719000000000001,72.203000000000003,0,0)

# Number format codes must exactly match the existing styles.xml formats
# (including the backslash-escaped parens) or the engine will mint a brand
# new numFmt/cellXf instead of reusing the existing style index.
$dateFormat = "yyyy\-mm\-dd"
$numberFormat = "0.00_);[Red]\(0.00\)"

# Row 150: 四方坪站
$ws.Cells.Item(150, 1).Value = $dateValue
$ws.Cells.Item(150, 1).NumberFormat = $dateFormat
$ws.Cells.Item(150, 2).Value = $station1
for ($i = 0; $i -lt $row150.Length; $i++) {
    $col = 3 + $i
    $cell = $ws.Cells.Item(150, $col)
    $cell.Value = $row150[$i]
    $cell.NumberFormat = $numberFormat
}

# Row 151: 高岭站
$ws.Cells.Item(151, 1).Value = $dateValue
$ws.Cells.Item(151, 1).NumberFormat = $dateFormat
$ws.Cells.Item(151, 2).Value = $station2
for ($i = 0; $i -lt $row151.Length; $i++) {
    $col = 3 + $i
    $cell = $ws.Cells.Item(151, $col)
    $cell.Value = $row151[$i]
    $cell.NumberFormat = $numberFormat
}

# Mirror the author's final selection recorded in the workbook view.
$ws.Range("D155").Select()
